# Auto-generated PowerShell Excel COM-interop edit script
# Applies the changes described by the diff to before.xlsx

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Schedule")
$ws2 = $wb.Worksheets.Item("Detailed")

# --- Sheet1 (Schedule) updates ---
$ws1.Range("E2").Value = 359.97305175
$ws1.Range("F2").Value = 7.935913839285715

$ws1.Range("B3").Value = 46045.125
$ws1.Range("C3").Value = 7
$ws1.Range("D3").Value = 26.46
$ws1.Range("E3").Value = 727.90174275
$ws1.Range("F3").Value = 27.50951408730159

# New row 4 in Schedule sheet
$ws1.Range("A4").Value = 46045.29166666666
$ws1.Range("A4").NumberFormat = $ws1.Range("A2").NumberFormat
$ws1.Range("B4").Value = 46045.66666666666
$ws1.Range("B4").NumberFormat = $ws1.Range("B2").NumberFormat
$ws1.Range("C4").Value = 9
$ws1.Range("D4").Value = 34.02
$ws1.Range("E4").Value = 10.32719024999999
$ws1.Range("F4").Value = 0.3035623236331567

# --- Sheet2 (Detailed) updates: Price (B) corrections for rows 25-46 ---
$priceUpdates = @{
    25 = -5.58973
    26 = -6.72804
    27 = -7.42351
    28 = -7.85989
    29 = -5.95281
    30 = -0.87314
    31 = -0.88011
    32 = 0.51
    33 = 7.66303
    34 = 3.61347
    35 = -5.38179
    36 = 0.02716
    37 = 1.60178
    38 = 0.29278
    39 = 10.26274
    40 = 22.81246
    41 = 50.55371
    42 = 57.31
    45 = 54.00669
    46 = 43.27957
}
foreach ($r in $priceUpdates.Keys) {
    $ws2.Cells.Item([int]$r, 2).Value = $priceUpdates[$r]
}

# --- Sheet2 (Detailed) updates: Type (C) corrections for rows 27-31 (forecast -> historical) ---
$typeUpdates = @{
    27 = "historical"
    28 = "historical"
    29 = "historical"
    30 = "historical"
    31 = "historical"
}
foreach ($r in $typeUpdates.Keys) {
    $ws2.Cells.Item([int]$r, 3).Value = $typeUpdates[$r]
}

# --- Sheet2 (Detailed): update row 49 Price, and append new rows 50-97 ---
$newRows = @(
    @(49, 46044.97916666666, 57.06, "forecast", 46044, "ON"),
    @(50, 46045, 56.98, "forecast", 46045, "ON"),
    @(51, 46045.02083333334, 56.98, "forecast", 46045, "ON"),
    @(52, 46045.04166666666, 52.17694, "forecast", 46045, "ON"),
    @(53, 46045.0625, 52.17549, "forecast", 46045, "ON"),
    @(54, 46045.08333333334, 41.99762, "forecast", 46045, "ON"),
    @(55, 46045.10416666666, 48.13511, "forecast", 46045, "ON"),
    @(56, 46045.125, 49.45214, "forecast", 46045, "OFF"),
    @(57, 46045.14583333334, 48.76499, "forecast", 46045, "OFF"),
    @(58, 46045.16666666666, 50.93505, "forecast", 46045, "OFF"),
    @(59, 46045.1875, 56.98, "forecast", 46045, "OFF"),
    @(60, 46045.20833333334, 56.98, "forecast", 46045, "OFF"),
    @(61, 46045.22916666666, 57.06, "forecast", 46045, "OFF"),
    @(62, 46045.25, 57.09, "forecast", 46045, "OFF"),
    @(63, 46045.27083333334, 57.06, "forecast", 46045, "OFF"),
    @(64, 46045.29166666666, 36.06, "forecast", 46045, "ON"),
    @(65, 46045.3125, 9.558339999999999, "forecast", 46045, "ON"),
    @(66, 46045.33333333334, 0.7, "forecast", 46045, "ON"),
    @(67, 46045.35416666666, 0.7, "forecast", 46045, "ON"),
    @(68, 46045.375, 0.7, "forecast", 46045, "ON"),
    @(69, 46045.39583333334, 0.02898, "forecast", 46045, "ON"),
    @(70, 46045.41666666666, -0.96199, "forecast", 46045, "ON"),
    @(71, 46045.4375, 0.00957, "forecast", 46045, "ON"),
    @(72, 46045.45833333334, 0.02837, "forecast", 46045, "ON"),
    @(73, 46045.47916666666, 0.7, "forecast", 46045, "ON"),
    @(74, 46045.5, 0.7, "forecast", 46045, "ON"),
    @(75, 46045.52083333334, -4.46499, "forecast", 46045, "ON"),
    @(76, 46045.54166666666, -2.54301, "forecast", 46045, "ON"),
    @(77, 46045.5625, -4.20359, "forecast", 46045, "ON"),
    @(78, 46045.58333333334, -5.2795, "forecast", 46045, "ON"),
    @(79, 46045.60416666666, -7.77674, "forecast", 46045, "ON"),
    @(80, 46045.625, -7.78627, "forecast", 46045, "ON"),
    @(81, 46045.64583333334, -5.57718, "forecast", 46045, "ON"),
    @(82, 46045.66666666666, -5.64276, "forecast", 46045, "OFF"),
    @(83, 46045.6875, -6, "forecast", 46045, "OFF"),
    @(84, 46045.70833333334, -6, "forecast", 46045, "OFF"),
    @(85, 46045.72916666666, -4.03567, "forecast", 46045, "OFF"),
    @(86, 46045.75, 8.598929999999999, "forecast", 46045, "OFF"),
    @(87, 46045.77083333334, 33.06657, "forecast", 46045, "OFF"),
    @(88, 46045.79166666666, 56.98, "forecast", 46045, "OFF"),
    @(89, 46045.8125, 57.09, "forecast", 46045, "OFF"),
    @(90, 46045.83333333334, 57.06, "forecast", 46045, "OFF"),
    @(91, 46045.85416666666, 57.06, "forecast", 46045, "OFF"),
    @(92, 46045.875, 54.83006, "forecast", 46045, "OFF"),
    @(93, 46045.89583333334, 54.55154, "forecast", 46045, "OFF"),
    @(94, 46045.91666666666, 48.7809, "forecast", 46045, "OFF"),
    @(95, 46045.9375, 57.06, "forecast", 46045, "OFF"),
    @(96, 46045.95833333334, 57.06, "forecast", 46045, "OFF"),
    @(97, 46045.97916666666, 56.98, "forecast", 46045, "OFF"),
)

foreach ($row in $newRows) {
    $r = [int]$row[0]
    $ws2.Cells.Item($r, 1).Value = $row[1]
    $ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
    $ws2.Cells.Item($r, 2).Value = $row[2]
    $ws2.Cells.Item($r, 3).Value = $row[3]
    $ws2.Cells.Item($r, 4).Value = $row[4]
    $ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
    $ws2.Cells.Item($r, 5).Value = $row[5]
}

"Done applying edits"